# Updates the cryptos worksheet values to reflect the latest scrape
# (GitHub Actions scheduled refresh), per the authoring commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "74.892.53"
$ws.Range("E2").Value = "  +1.37%  "
$ws.Range("D3").Value = "2.817.64"
$ws.Range("E3").Value = "  +7.23%  "
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("D5").Value = "'187.06"
$ws.Range("E5").Value = "  +1.06%  "
$ws.Range("D6").Value = "'594.12"
$ws.Range("E6").Value = "  +1.89%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("D8").Value = "'0.549"
$ws.Range("E8").Value = "  +2.80%  "
$ws.Range("E9").Value = "  -5.07%  "
$ws.Range("D10").Value = "2.814.66"
$ws.Range("E10").Value = "  +7.17%  "
$ws.Range("E12").Value = "  +3.21%  "
$ws.Range("D13").Value = "'4.86"
$ws.Range("E13").Value = "  +1.89%  "
$ws.Range("D14").Value = "3.331.36"
$ws.Range("E14").Value = "  +7.18%  "
$ws.Range("D15").Value = "74.815.79"
$ws.Range("E15").Value = "  +1.49%  "
$ws.Range("E16").Value = "  -0.91%  "
$ws.Range("D17").Value = "'26.76"
$ws.Range("E17").Value = "  +2.10%  "
$ws.Range("D18").Value = "2.811.22"
$ws.Range("E18").Value = "  +7.09%  "
$ws.Range("E19").Value = "  -1.49%  "
$ws.Range("D20").Value = "'12.26"
$ws.Range("E20").Value = "  +3.30%  "
$ws.Range("D21").Value = "'376.28"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("E22").Value = "  -2.84%  "
$ws.Range("D23").Value = "'4.06"
$ws.Range("E23").Value = "  -1.03%  "
$ws.Range("E24").Value = "  -0.21%  "
$ws.Range("E25").Value = "  -0.13%  "
$ws.Range("D26").Value = "'70.57"
$ws.Range("E26").Value = "  +1.00%  "
$ws.Range("D27").Value = "2.958.48"
$ws.Range("E27").Value = "  +7.36%  "
$ws.Range("D28").Value = "'4.15"
$ws.Range("E28").Value = "  +0.03%  "
$ws.Range("E29").Value = "  +2.72%  "
$ws.Range("E30").Value = "  +9.77%  "
$ws.Range("E31").Value = "  +0.02%  "
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "'512.73"
$ws.Range("E32").Value = "  -1.99%  "
$ws.Range("B33").Value = "Fetch.AI"
$ws.Range("C33").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D33").Value = "'1.39"
$ws.Range("E33").Value = "  -0.86%  "
$ws.Range("D34").Value = "'7.69"
$ws.Range("E34").Value = "  +0.13%  "
$ws.Range("E35").Value = "  +2.36%  "
$ws.Range("E36").Value = "  -0.01%  "
$ws.Range("D37").Value = "'162.25"
$ws.Range("E37").Value = "  +0.78%  "
$ws.Range("D38").Value = "'19.91"
$ws.Range("E38").Value = "  +3.80%  "
$ws.Range("E39").Value = "  -1.66%  "
$ws.Range("E40").Value = "  +0.59%  "
$ws.Range("D41").Value = "'185.58"
$ws.Range("E41").Value = "  +14.88%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("E43").Value = "  +2.92%  "
$ws.Range("E44").Value = "  +1.23%  "
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("E46").Value = "  +1.19%  "
$ws.Range("D47").Value = "'39.95"
$ws.Range("E47").Value = "  +2.50%  "
$ws.Range("E48").Value = "  -2.03%  "
$ws.Range("E49").Value = "  -0.08%  "
$ws.Range("D50").Value = "'0.570"
$ws.Range("E50").Value = "  +7.78%  "
$ws.Range("D51").Value = "'3.70"
$ws.Range("E51").Value = "  +1.91%  "
